# Apply scheduled-runner updates to market/profit data cells across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row92
$ws.Cells.Item(92, 8).Value = 250000000
$ws.Cells.Item(92, 9).Value = 250000000
$ws.Cells.Item(92, 11).Value = 250000000
$ws.Cells.Item(92, 13).Value = -249998752

# ALC!row100
$ws.Cells.Item(100, 8).Value = 27779344
$ws.Cells.Item(100, 9).Value = 1558.3334
$ws.Cells.Item(100, 10).Value = 83334920
$ws.Cells.Item(100, 11).Value = 1558.3334
$ws.Cells.Item(100, 12).Value = 83334920
$ws.Cells.Item(100, 13).Value = -1017.3334
$ws.Cells.Item(100, 14).Value = -83336002

# ALC!row103
$ws.Cells.Item(103, 8).Value = 7692921
$ws.Cells.Item(103, 9).Value = 467.41666
$ws.Cells.Item(103, 10).Value = 14286453
$ws.Cells.Item(103, 11).Value = 1402.24998
$ws.Cells.Item(103, 12).Value = 42859359
$ws.Cells.Item(103, 13).Value = -816.2499800000001
$ws.Cells.Item(103, 14).Value = -42860531

$ws = $wb.Worksheets.Item("ARM")
# ARM!row3
$ws.Cells.Item(3, 8).Value = 4502.1055
$ws.Cells.Item(3, 9).Value = 2650
$ws.Cells.Item(3, 10).Value = 4720
$ws.Cells.Item(3, 11).Value = 2650
$ws.Cells.Item(3, 12).Value = 4720
$ws.Cells.Item(3, 13).Value = -2535
$ws.Cells.Item(3, 14).Value = -4950

# ARM!row8
$ws.Cells.Item(8, 8).Value = 6006235
$ws.Cells.Item(8, 9).Value = 10000400
$ws.Cells.Item(8, 10).Value = 14988
$ws.Cells.Item(8, 11).Value = 10000400
$ws.Cells.Item(8, 12).Value = 14988
$ws.Cells.Item(8, 13).Value = -10000256
$ws.Cells.Item(8, 14).Value = -15276

# ARM!row10
$ws.Cells.Item(10, 8).Value = 3000
$ws.Cells.Item(10, 10).Value = 3000
$ws.Cells.Item(10, 12).Value = 3000
$ws.Cells.Item(10, 14).Value = -3340

# ARM!row11
$ws.Cells.Item(11, 8).Value = 400
$ws.Cells.Item(11, 9).Value = 400
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 400
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -256
$ws.Cells.Item(11, 14).ClearContents()

# ARM!row13
$ws.Cells.Item(13, 8).Value = 6669400
$ws.Cells.Item(13, 9).Value = 10000100
$ws.Cells.Item(13, 10).Value = 8000
$ws.Cells.Item(13, 11).Value = 10000100
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = -9999956
$ws.Cells.Item(13, 14).Value = -8288

# ARM!row63
$ws.Cells.Item(63, 8).Value = 2613.64
$ws.Cells.Item(63, 9).Value = 1521.3125
$ws.Cells.Item(63, 10).Value = 4555.5557
$ws.Cells.Item(63, 11).Value = 1521.3125
$ws.Cells.Item(63, 12).Value = 4555.5557
$ws.Cells.Item(63, 13).Value = -835.3125
$ws.Cells.Item(63, 14).Value = -5927.5557

# ARM!row66
$ws.Cells.Item(66, 8).Value = 2613.64
$ws.Cells.Item(66, 9).Value = 1521.3125
$ws.Cells.Item(66, 10).Value = 4555.5557
$ws.Cells.Item(66, 11).Value = 7606.5625
$ws.Cells.Item(66, 12).Value = 22777.7785
$ws.Cells.Item(66, 13).Value = -4174.5625
$ws.Cells.Item(66, 14).Value = -29641.7785

$ws = $wb.Worksheets.Item("BSM")
# BSM!row126
$ws.Cells.Item(126, 8).Value = 31804
$ws.Cells.Item(126, 10).Value = 31804
$ws.Cells.Item(126, 12).Value = 31804
$ws.Cells.Item(126, 14).Value = -41684

$ws = $wb.Worksheets.Item("CRP")
# CRP!row7
$ws.Cells.Item(7, 8).Value = 31.428572
$ws.Cells.Item(7, 9).Value = 32.727272
$ws.Cells.Item(7, 10).Value = 26.666666
$ws.Cells.Item(7, 11).Value = 32.727272
$ws.Cells.Item(7, 12).Value = 26.666666
$ws.Cells.Item(7, 13).Value = 80.272728
$ws.Cells.Item(7, 14).Value = -252.666666

# CRP!row118
$ws.Cells.Item(118, 8).Value = 30000
$ws.Cells.Item(118, 10).Value = 30000
$ws.Cells.Item(118, 12).Value = 30000
$ws.Cells.Item(118, 14).Value = -33314

$ws = $wb.Worksheets.Item("CUL")
# CUL!row2
$ws.Cells.Item(2, 8).Value = 1556416.6
$ws.Cells.Item(2, 9).Value = 71447
$ws.Cells.Item(2, 11).Value = 428682
$ws.Cells.Item(2, 13).Value = -428569

# CUL!row92
$ws.Cells.Item(92, 8).Value = 915.3333
$ws.Cells.Item(92, 9).Value = 691.25
$ws.Cells.Item(92, 10).Value = 1171.4286
$ws.Cells.Item(92, 11).Value = 2073.75
$ws.Cells.Item(92, 12).Value = 3514.2858
$ws.Cells.Item(92, 13).Value = -825.75
$ws.Cells.Item(92, 14).Value = -6010.2858

# CUL!row114
$ws.Cells.Item(114, 8).Value = 613.53845
$ws.Cells.Item(114, 9).Value = 584.5
$ws.Cells.Item(114, 10).Value = 710.3333
$ws.Cells.Item(114, 11).Value = 1753.5
$ws.Cells.Item(114, 12).Value = 2130.9999
$ws.Cells.Item(114, 13).Value = 1500.5
$ws.Cells.Item(114, 14).Value = -8638.999899999999

$ws = $wb.Worksheets.Item("GSM")
# GSM!row55
$ws.Cells.Item(55, 8).Value = 26333
$ws.Cells.Item(55, 10).Value = 26333
$ws.Cells.Item(55, 12).Value = 26333
$ws.Cells.Item(55, 14).Value = -26987

# GSM!row70
$ws.Cells.Item(70, 8).Value = 38967.832
$ws.Cells.Item(70, 9).Value = 45832.32
$ws.Cells.Item(70, 10).Value = 4645.4
$ws.Cells.Item(70, 11).Value = 45832.32
$ws.Cells.Item(70, 12).Value = 4645.4
$ws.Cells.Item(70, 13).Value = -45562.32
$ws.Cells.Item(70, 14).Value = -5185.4

# GSM!row73
$ws.Cells.Item(73, 8).Value = 38967.832
$ws.Cells.Item(73, 9).Value = 45832.32
$ws.Cells.Item(73, 10).Value = 4645.4
$ws.Cells.Item(73, 11).Value = 45832.32
$ws.Cells.Item(73, 12).Value = 4645.4
$ws.Cells.Item(73, 13).Value = -44896.32
$ws.Cells.Item(73, 14).Value = -6517.4

# GSM!row126
$ws.Cells.Item(126, 8).Value = 10919.129
$ws.Cells.Item(126, 9).Value = 3146.25
$ws.Cells.Item(126, 10).Value = 13622.739
$ws.Cells.Item(126, 11).Value = 9438.75
$ws.Cells.Item(126, 12).Value = 40868.217
$ws.Cells.Item(126, 13).Value = -6968.75
$ws.Cells.Item(126, 14).Value = -45808.217

# GSM!row132
$ws.Cells.Item(132, 8).Value = 2083.2122
$ws.Cells.Item(132, 9).Value = 1861.037
$ws.Cells.Item(132, 11).Value = 5583.111
$ws.Cells.Item(132, 13).Value = -3053.111

# GSM!row133
$ws.Cells.Item(133, 8).Value = 38610
$ws.Cells.Item(133, 10).Value = 38610
$ws.Cells.Item(133, 12).Value = 38610
$ws.Cells.Item(133, 14).Value = -48730

$ws = $wb.Worksheets.Item("LTW")
# LTW!row40
$ws.Cells.Item(40, 8).Value = 2648.111
$ws.Cells.Item(40, 9).Value = 2251
$ws.Cells.Item(40, 10).Value = 2965.8
$ws.Cells.Item(40, 11).Value = 2251
$ws.Cells.Item(40, 12).Value = 2965.8
$ws.Cells.Item(40, 13).Value = -2115
$ws.Cells.Item(40, 14).Value = -3237.8

# LTW!row46
$ws.Cells.Item(46, 8).Value = 1616.8286
$ws.Cells.Item(46, 9).Value = 1475.6
$ws.Cells.Item(46, 10).Value = 1969.9
$ws.Cells.Item(46, 11).Value = 1475.6
$ws.Cells.Item(46, 12).Value = 1969.9
$ws.Cells.Item(46, 13).Value = -1287.6
$ws.Cells.Item(46, 14).Value = -2345.9

$ws = $wb.Worksheets.Item("WVR")
# WVR!row62
$ws.Cells.Item(62, 8).Value = 8189
$ws.Cells.Item(62, 9).Value = 4196.6665
$ws.Cells.Item(62, 10).Value = 9900
$ws.Cells.Item(62, 11).Value = 4196.6665
$ws.Cells.Item(62, 12).Value = 9900
$ws.Cells.Item(62, 13).Value = -3572.6665
$ws.Cells.Item(62, 14).Value = -11148

# WVR!row65
$ws.Cells.Item(65, 8).Value = 8189
$ws.Cells.Item(65, 9).Value = 4196.6665
$ws.Cells.Item(65, 10).Value = 9900
$ws.Cells.Item(65, 11).Value = 20983.3325
$ws.Cells.Item(65, 12).Value = 49500
$ws.Cells.Item(65, 13).Value = -17863.3325
$ws.Cells.Item(65, 14).Value = -55740

# WVR!row107
$ws.Cells.Item(107, 8).Value = 3814.4285
$ws.Cells.Item(107, 9).Value = 2951
$ws.Cells.Item(107, 10).Value = 4462
$ws.Cells.Item(107, 11).Value = 8853
$ws.Cells.Item(107, 12).Value = 13386
$ws.Cells.Item(107, 13).Value = -6933
$ws.Cells.Item(107, 14).Value = -17226
